# Update cryptos list with refreshed prices and 1h volume percentages
# (as produced by the scheduled GitHub Actions scraper run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.560.71"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "'2.270.59"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").Value = "'119.26"
$ws.Range("E5").Value = "  +4.95%  "
$ws.Range("D6").Value = "'265.86"
$ws.Range("E6").Value = "  -0.42%  "
$ws.Range("D7").Value = "'0.647"
$ws.Range("E7").Value = "  +3.64%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D9").Value = "'0.621"
$ws.Range("E9").Value = "  +1.84%  "
$ws.Range("D10").Value = "'47.69"
$ws.Range("E10").Value = "  -1.16%  "
$ws.Range("D11").Value = "'0.0944"
$ws.Range("E11").Value = "  +0.96%  "
$ws.Range("D12").Value = "'9.21"
$ws.Range("E12").Value = "  +4.52%  "
$ws.Range("D14").Value = "'15.45"
$ws.Range("E14").Value = "  -1.40%  "
$ws.Range("D15").Value = "'0.905"
$ws.Range("E15").Value = "  +4.16%  "
$ws.Range("D16").Value = "'2.610.42"
$ws.Range("E16").Value = "  -0.52%  "
$ws.Range("D17").Value = "'2.262.62"
$ws.Range("E17").Value = "  -0.70%  "
$ws.Range("D18").Value = "'43.539.44"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("E19").Value = "  +1.36%  "
$ws.Range("D20").Value = "'6.86"
$ws.Range("E20").Value = "  -2.19%  "
$ws.Range("D21").Value = "'72.13"
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("D22").Value = "'2.41"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").Value = "'235.93"
$ws.Range("E23").Value = "  +1.65%  "
$ws.Range("D24").Value = "'9.54"
$ws.Range("E24").Value = "  -4.13%  "
$ws.Range("D25").Value = "'2.88"
$ws.Range("E25").Value = "  +0.24%  "
$ws.Range("D26").Value = "'12.05"
$ws.Range("E26").Value = "  +5.01%  "
$ws.Range("E27").Value = "  +1.79%  "
$ws.Range("D28").Value = "'41.82"
$ws.Range("E28").Value = "  +1.81%  "
$ws.Range("D29").Value = "'3.38"
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("D31").Value = "'172.21"
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("D32").Value = "'21.60"
$ws.Range("E32").Value = "  +0.79%  "
$ws.Range("D33").Value = "'0.0917"
$ws.Range("E33").Value = "  +0.63%  "
$ws.Range("D34").Value = "'5.73"
$ws.Range("E34").Value = "  +1.69%  "
$ws.Range("D35").Value = "'0.131"
$ws.Range("E35").Value = "  +2.83%  "
$ws.Range("E36").Value = "  +9.63%  "
$ws.Range("D37").Value = "'4.19"
$ws.Range("E37").Value = "  +12.09%  "
$ws.Range("D38").Value = "'4.56"
$ws.Range("E38").Value = "  -0.74%  "
$ws.Range("E39").Value = "  +1.43%  "
$ws.Range("D40").Value = "'2.56"
$ws.Range("E40").Value = "  +5.85%  "
$ws.Range("D41").Value = "'73.98"
$ws.Range("E41").Value = "  -0.87%  "
$ws.Range("D42").Value = "'13.73"
$ws.Range("E42").Value = "  -4.23%  "
$ws.Range("D43").Value = "'0.237"
$ws.Range("E43").Value = "  +0.38%  "
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("E45").Value = "  -0.43%  "
$ws.Range("E46").Value = "  -6.38%  "
$ws.Range("D47").Value = "'73.85"
$ws.Range("E47").Value = "  +41.32%  "
$ws.Range("E48").Value = "  +1.73%  "
$ws.Range("D49").Value = "'8.55"
$ws.Range("E49").Value = "  -1.14%  "
$ws.Range("D50").Value = "'0.100"
$ws.Range("E50").Value = "  +0.53%  "
$ws.Range("D51").Value = "'101.96"
$ws.Range("E51").Value = "  +0.35%  "

# Strip the quote-prefix formatting introduced by the apostrophe above so
# the price cells keep their original (General) style/number format.
$ws.Range("D2:D51").ClearFormats()
